# Update "想去人数" (interested-count) figures that changed between scrapes.
# Sheet 展览 (Exhibitions), 演出 (Shows) and 全部类型 (All types) each carry
# their own copy of the same events, so each copy's F-column needs updating.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2  = 153
        3  = 1806
        13 = 169
        15 = 124
        18 = 5139
        20 = 842
        21 = 118
        22 = 2286
        24 = 31
        25 = 2134
    }
    "演出" = @{
        2 = 86
    }
    "全部类型" = @{
        2  = 153
        3  = 1806
        13 = 169
        15 = 124
        18 = 5139
        19 = 86
        22 = 842
        23 = 118
        24 = 2286
        27 = 31
        28 = 2134
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
